$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.325.35"
$ws.Range("D3").Value = "2.505.87"
$ws.Range("E3").Value = "  +2.27%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "324.35"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").Value = "109.17"
$ws.Range("E6").Value = "  +4.19%  "
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.537"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").Value = "39.12"
$ws.Range("E10").Value = "  +8.73%  "
$ws.Range("E12").Value = "  +0.55%  "
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").Value = "7.20"
$ws.Range("E14").Value = "  +1.94%  "
$ws.Range("D15").Value = "2.897.36"
$ws.Range("E15").Value = "  +2.27%  "
$ws.Range("D16").Value = "2.502.12"
$ws.Range("E16").Value = "  +2.54%  "
$ws.Range("E17").Value = "  +1.72%  "
$ws.Range("D18").Value = "47.247.86"
$ws.Range("E18").Value = "  +2.63%  "
$ws.Range("D19").Value = "12.85"
$ws.Range("E19").Value = "  +1.54%  "
$ws.Range("D21").Value = "0.0₃0942"
$ws.Range("E21").Value = "  +1.05%  "
$ws.Range("D22").Value = "2.72"
$ws.Range("E22").Value = "  +13.65%  "
$ws.Range("D23").Value = "70.52"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("D24").Value = "247.44"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("E25").Value = "  +3.45%  "
$ws.Range("D26").Value = "26.04"
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D29").Value = "10.07"
$ws.Range("E29").Value = "  +3.90%  "
$ws.Range("D30").Value = "35.58"
$ws.Range("E30").Value = "  +3.88%  "
$ws.Range("E31").Value = "  +8.64%  "
$ws.Range("D32").Value = "49.85"
$ws.Range("E32").Value = "  +0.99%  "
$ws.Range("D33").Value = "20.11"
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").Value = "5.44"
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("D35").Value = "0.0792"
$ws.Range("E35").Value = "  +3.94%  "
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("E37").Value = "  +5.03%  "
$ws.Range("E38").Value = "  +3.85%  "
$ws.Range("D39").Value = "3.01"
$ws.Range("E39").Value = "  +1.60%  "
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("D42").Value = "121.32"
$ws.Range("E42").Value = "  -5.21%  "
$ws.Range("D43").Value = "21.24"
$ws.Range("E43").Value = "  +2.11%  "
$ws.Range("E44").Value = "  +2.15%  "
$ws.Range("D45").Value = "1.995.80"
$ws.Range("E45").Value = "  +1.26%  "
$ws.Range("E46").Value = "  +3.96%  "
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "1.78"
$ws.Range("E48").Value = "  -4.24%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "9.09"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  +3.64%  "
$ws.Range("E51").Value = "  +4.42%  "
